$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Dataøving 5" (row 16, B) moves/merges into row 17, combined with
# "Forberedelse til eksamen" (row 17, B) which in turn moves up to row 16.
$ws.Range("B16").Value = "Forberedelse til eksamen"
$ws.Range("B17").Value = "Forberedelse til eksamen/Dataøving 5"

$ws.Range("B17").Select()
